$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

# --- Normalize row 13: the empty-looking cells (B..K, M) actually held an
# empty string; the new data makes them hold the literal text "nan", same
# as every other data row on this sheet. L13/N13/O13 already have values
# and stay untouched.
$row13NanCols = 2,3,4,5,6,7,8,9,10,11,13
foreach ($c in $row13NanCols) {
    $ws.Cells.Item(13, $c).Value = "nan"
}

# --- Append new row 14 for the new Card21 service event.
# Card/event id ("21") is stored as text like every other row on this sheet
# (not a number), so format as Text before writing it, then drop the
# resulting style back to the sheet's default (Normal) so only the value
# - not an extra number format - is new.
$ws.Cells.Item(14, 1).NumberFormat = "@"
$ws.Cells.Item(14, 1).Value = "21"
$ws.Cells.Item(14, 1).Style = "Normal"

$ws.Cells.Item(14, 12).Value = "28\8\2024"
$ws.Cells.Item(14, 14).Value = "تم عمل معايره للمكنه steeing"
$ws.Cells.Item(14, 15).Value = "الخبير"

Write-Output "Card21 row14 added"
